$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 142, shifting existing rows 142-167 down to 143-168
$ws.Rows("142:142").Insert()

# Populate the new row 142 with the new weekly data entry
$ws.Range("A142").Value = 7
$ws.Range("B142").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C142").Value = "Ñuble"
$ws.Range("D142").Value = 44474
$ws.Range("D142").NumberFormat = $ws.Range("D143").NumberFormat
$ws.Range("E142").Value = 16
$ws.Range("F142").Value = 100112023
$ws.Range("G142").Value = "Brócoli"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 120
$ws.Range("K142").Value = 650
$ws.Range("L142").Value = 700
$ws.Range("M142").Value = 675
$ws.Range("N142").Value = "`$/unidad"
$ws.Range("O142").Value = "Región del Maule"
$ws.Range("P142").Value = 675
$ws.Range("Q142").Value = 1
$ws.Range("R142").Value = "Hortaliza"
